$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the header row ("Name","Size","Rank","Tuition") so the data shifts up
# one row (row 2 -> row 1, ... row 12 -> row 11).
$ws.Rows.Item(1).Delete() | Out-Null

# Populate the (previously empty) Rank column (C) with each school's rank.
$ranks = @(9, 9, 10, 207, 1, 1, 30, 14, 27, 198, 20)

for ($i = 0; $i -lt $ranks.Length; $i++) {
    $row = $i + 1
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = $ranks[$i]
    # Matches the formatting already used elsewhere in the sheet (reuses the
    # existing-but-unused style record instead of creating a new one).
    $cell.Font.ColorIndex = -4105
}

# Restore the active selection recorded for this sheet.
$ws.Range("C14").Select() | Out-Null
